$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$values = @{
    2 = 132.8887902256449
    3 = 13.11415068141123
    4 = 12.6005575634989
    5 = 16.73368674978447
    6 = 22.84727367116519
    7 = 7.225356000066316
    8 = 6.667566441740168
    9 = 21.82999353405398
    10 = 36.28984962530672
    11 = 10.22995194061277
    12 = 2.245636377578088
    13 = 6.075486541713641
    14 = 1.455434503282711
    15 = 2.547391384541415
    16 = 18.18073916954042
    17 = 19.07990812968366
    18 = 18.4817321689104
    19 = 6.533453105312023
    20 = 27.01240717637111
    21 = 69.72217077558193
    22 = 12.82915984150142
    23 = 2.39117982790636
    24 = 23.16274481480867
    25 = 6.669852464520192
    26 = 12.70952464934677
    27 = 28.52194421878108
    28 = 4.947715303567875
    29 = 12.59750953312554
    30 = 2.385083767159625
    31 = 2.458998503713774
    32 = 4.342681274454539
    33 = 5.122215042443132
    34 = 97.55068808442203
    35 = 8.600017698454803
    36 = 22.708588289177
    37 = 4.235238203793354
    38 = 9.923624888089394
    39 = 8.917774864878306
    40 = 7.73285305723191
    41 = 5.811069906824059
    42 = 264.62
}

foreach ($row in $values.Keys) {
    $ws.Cells.Item($row, 3).Value = $values[$row]
}

